$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 94
$ws.Range("F8").Value = 9487
$ws.Range("F11").Value = 691
$ws.Range("F12").Value = 1959
$ws.Range("F14").Value = 898
$ws.Range("F15").Value = 2608
$ws.Range("F16").Value = 127
$ws.Range("F17").Value = 3911
$ws.Range("F18").Value = 307
$ws.Range("F19").Value = 139
$ws.Range("F20").Value = 125
$ws.Range("F25").Value = 70
$ws.Range("F26").Value = 259
$ws.Range("F27").Value = 561
$ws.Range("C28").Value = "北京·TCS卡牌嘉年华"
$ws.Range("F29").Value = 2139
$ws.Range("F30").Value = 1094
$ws.Range("F33").Value = 4317
$ws.Range("F34").Value = 68
$ws.Range("F35").Value = 175
$ws.Range("F36").Value = 339
$ws.Range("F37").Value = 150

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value = 0
$ws.Range("F6").Value = 20

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 188
$ws.Range("F3").Value = 974

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 188
$ws.Range("F4").Value = 974
$ws.Range("F5").Value = 2
$ws.Range("G6").Value = 0
$ws.Range("F7").Value = 94
$ws.Range("F11").Value = 9487
$ws.Range("F14").Value = 691
$ws.Range("F15").Value = 1959
$ws.Range("F17").Value = 898
$ws.Range("F19").Value = 2608
$ws.Range("F20").Value = 127
$ws.Range("F21").Value = 3911
$ws.Range("F22").Value = 307
$ws.Range("F23").Value = 139
$ws.Range("F24").Value = 125
$ws.Range("F30").Value = 70
$ws.Range("F31").Value = 259
$ws.Range("F32").Value = 561
$ws.Range("C33").Value = "北京·TCS卡牌嘉年华"
$ws.Range("F34").Value = 2139
$ws.Range("F35").Value = 1094
$ws.Range("F38").Value = 4317
$ws.Range("F39").Value = 68
$ws.Range("F40").Value = 175
$ws.Range("F41").Value = 339
$ws.Range("F42").Value = 150
$ws.Range("F44").Value = 20
